$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width
$ws.Columns.Item(1).ColumnWidth = 22.88671875

# Row 12: J12 = AVERAGE(J2:J11)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Row 14
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

# Row 15
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# Row 16
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# Row 17
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$ws.Range("B14:B17").Font.Bold = $true
$ws.Range("B14,B16").Font.Size = 12
$ws.Range("B15,B17").Font.Size = 12
$ws.Range("B15,B17").VerticalAlignment = -4108
